# Rename header columns: "_old" -> "_FV2310", "_new" -> "_FV2404"
# and turn the header row range into a real Excel Table, with a frozen
# top row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the header row labels -------------------------------------
# "<name>_old" -> "<name>_FV2310", "<name>_new" -> "<name>_FV2404"
$newHeaders = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newHeaders[$i]
}

# --- 2) Turn the used range into an Excel Table ---------------------------
$tableRange = $ws.Range("A1:U79")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# --- 3) Freeze the header row ----------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
